$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unitario")

# Insert a new row above the current row 27 (shifts rows 27:61 down to 28:62)
$ws.Rows.Item(27).Insert()

# Populate the new row with the água (water) related service entry
$ws.Cells.Item(27, 1).Value = "280000"
$ws.Cells.Item(27, 2).Value = "PASSADO RAMAL DE AGUA PARA NOVA REDE"
$ws.Cells.Item(27, 3).Value = "LigacaoAgua"

# Match the row height used by the rest of the table
$ws.Rows.Item(27).RowHeight = 18.75

# Nudge alignment formatting so the new cells pick up their own style
# (matches the border/alignment xf the author's row ended up with)
$ws.Range("B27:C27").Orientation = 0

# Make this sheet ("unitario") the active / selected tab, with the view and
# selection the author left it in
$ws.Select()
$ws.Range("E25").Select()
